$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.153.35'
$ws.Range("D3").Value = '1.824.95'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '310.90'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '0.4960'
$ws.Range("E7").Value = '  -3.52%  '
$ws.Range("D8").Value = '0.3925'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").Value = '0.09867'
$ws.Range("E9").Value = '  +25.52%  '
$ws.Range("D10").Value = '1.108'
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("D11").Value = '41.29'
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '6.450'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '1.001'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '1.823.02'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = '7.311'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '0.00001143'
$ws.Range("E17").Value = '  +5.52%  '
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '0.06662'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '17.24'
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").Value = '5.991'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").Value = '28.202.83'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").Value = '11.38'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").Value = '2.247'
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").Value = '158.71'
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("D28").Value = '2.032.73'
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("D29").Value = '2.423'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = '127.15'
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("D32").Value = '1.040'
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("D33").Value = '5.611'
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").Value = '3.600'
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("D35").Value = '0.06749'
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("D36").Value = '9.018'
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("D37").Value = '0.02344'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("D39").Value = '4.985'
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").Value = '11.40'
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("D41").Value = '0.6229'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = '1.181'
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '13.26'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '0.5944'
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").Value = '3.709'
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("D48").Value = '124.20'
$ws.Range("E48").Value = '  -1.22%  '
$ws.Range("D49").Value = '1.951'
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("E50").Value = '  -3.04%  '
$ws.Range("D51").Value = '0.06795'
$ws.Range("E51").Value = '  -0.84%  '
